# "report and slight math fix"
#
# 1. Add a new header/caption row above the "(seconds)" report table explaining
#    that the values below are expressed in seconds, boxed with a thin border.
# 2. Rename the tamingoftheshrew / midsummersnightsdream column headers in that
#    report table so they also carry the "(seconds)" suffix (matching the
#    twelfthnight column, which already had it).
# 3. Math fix: C19 (=C9/1000000000) was showing a misleading 0 because C9 is
#    never populated - clear the bogus formula instead of displaying 0.
# 4. Leave the sheet positioned/sized the way the author left it after saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New labeled/boxed row just above the seconds-based report table.
$ws.Range("A11").Value = "Below is the time in seconds"
$ws.Range("A11").Borders.LineStyle = 1

# 2. Update the two column headers that were missing the "(seconds)" suffix.
$ws.Range("C12").Value = "Indexing tamingoftheshrew(seconds)"
$ws.Range("D12").Value = "Indexing midsummersnightsdream(seconds)"

# 3. Math fix - stop dividing an always-empty C9 into a fake 0.
$ws.Range("C19").ClearContents()

# 4. Keep row 23 present (extends the sheet's used range) and leave the
#    selection where the author ended up (E24), matching the saved file.
$ws.Range("A23").Font.Bold = $false
$ws.Range("E24").Select()
